$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated Ligand/Receptor-expressing cell counts (1 -> 3) and all dependent
# NATMI-derived statistics recomputed per Dr Hou's advice, row by row.
$updates = @{
    "E2" = 3; "G2" = 6.980814333333332; "H2" = 20.942443; "I2" = 0.2573350203399358; "J2" = 0.2573350203399358; "K2" = 3; "M2" = 490.031855; "N2" = 1470.095565; "O2" = 0.6686419015677429; "P2" = 0.6686419015677431; "Q2" = 3420.821397173921; "R2" = 30787.39257456529; "S2" = 0.1720649773400685; "T2" = 0.1720649773400685
    "E3" = 3; "G3" = 6.980814333333332; "H3" = 20.942443; "I3" = 0.2573350203399358; "J3" = 0.2573350203399358; "K3" = 3; "M3" = 62.79306433333334; "N3" = 188.379193; "O3" = 0.0856802950924601; "P3" = 0.08568029509246011; "Q3" = 438.3467235320554; "R3" = 3945.120511788499; "S3" = 0.02204854048034992; "T3" = 0.02204854048034992
    "E4" = 3; "G4" = 6.980814333333332; "H4" = 20.942443; "I4" = 0.2573350203399358; "J4" = 0.2573350203399358; "K4" = 3; "M4" = 0.5977846666666666; "N4" = 1.793354; "O4" = 0.0008156691696053909; "P4" = 0.000815669169605391; "Q4" = 4.173023769313555; "R4" = 37.55721392382199; "S4" = 0.0002099002423510618; "T4" = 0.0002099002423510618
    "E5" = 3; "G5" = 6.980814333333332; "H5" = 20.942443; "I5" = 0.2573350203399358; "J5" = 0.2573350203399358; "K5" = 3; "M5" = 179.453674; "N5" = 538.361022; "O5" = 0.2448621341701915; "P5" = 0.2448621341701916; "Q5" = 1252.732779628527; "R5" = 11274.59501665675; "S5" = 0.06301160227716632; "T5" = 0.06301160227716633
    "E6" = 3; "G6" = 6.375361666666667; "H6" = 19.126085; "I6" = 0.2350161092714131; "J6" = 0.2350161092714131; "K6" = 3; "M6" = 490.031855; "N6" = 1470.095565; "O6" = 0.6686419015677429; "P6" = 0.6686419015677431; "Q6" = 3124.130303812558; "R6" = 28117.17273431303; "S6" = 0.1571416182022901; "T6" = 0.1571416182022902
    "E7" = 3; "G7" = 6.375361666666667; "H7" = 19.126085; "I7" = 0.2350161092714131; "J7" = 0.2350161092714131; "K7" = 3; "M7" = 62.79306433333334; "N7" = 188.379193; "O7" = 0.0856802950924601; "P7" = 0.08568029509246011; "Q7" = 400.3284952832673; "R7" = 3602.956457549405; "S7" = 0.02013624959385652; "T7" = 0.02013624959385653
    "E8" = 3; "G8" = 6.375361666666667; "H8" = 19.126085; "I8" = 0.2350161092714131; "J8" = 0.2350161092714131; "K8" = 3; "M8" = 0.5977846666666666; "N8" = 1.793354; "O8" = 0.0008156691696053909; "P8" = 0.000815669169605391; "Q8" = 3.811093448787778; "R8" = 34.29984103909; "S8" = 0.0001916953946933034; "T8" = 0.0001916953946933034
    "E9" = 3; "G9" = 6.375361666666667; "H9" = 19.126085; "I9" = 0.2350161092714131; "J9" = 0.2350161092714131; "K9" = 3; "M9" = 179.453674; "N9" = 538.361022; "O9" = 0.2448621341701915; "P9" = 0.2448621341701916; "Q9" = 1144.082074162097; "R9" = 10296.73866745887; "S9" = 0.05754654608057316; "T9" = 0.05754654608057316
    "E10" = 3; "G10" = 5.973131; "H10" = 17.919393; "I10" = 0.2201886075150976; "J10" = 0.2201886075150976; "K10" = 3; "M10" = 490.031855; "N10" = 1470.095565; "O10" = 0.6686419015677429; "P10" = 0.6686419015677431; "Q10" = 2927.024464088005; "R10" = 26343.22017679205; "S10" = 0.1472273292324483; "T10" = 0.1472273292324483
    "E11" = 3; "G11" = 5.973131; "H11" = 17.919393; "I11" = 0.2201886075150976; "J11" = 0.2201886075150976; "K11" = 3; "M11" = 62.79306433333334; "N11" = 188.379193; "O11" = 0.0856802950924601; "P11" = 0.08568029509246011; "Q11" = 375.0711991544277; "R11" = 3375.640792389849; "S11" = 0.01886582486789144; "T11" = 0.01886582486789144
    "E12" = 3; "G12" = 5.973131; "H12" = 17.919393; "I12" = 0.2201886075150976; "J12" = 0.2201886075150976; "K12" = 3; "M12" = 0.5977846666666666; "N12" = 1.793354; "O12" = 0.0008156691696053909; "P12" = 0.000815669169605391; "Q12" = 3.570646123791333; "R12" = 32.135815114122; "S12" = 0.000179601058648407; "T12" = 0.000179601058648407
    "E13" = 3; "G13" = 5.973131; "H13" = 17.919393; "I13" = 0.2201886075150976; "J13" = 0.2201886075150976; "K13" = 3; "M13" = 179.453674; "N13" = 538.361022; "O13" = 0.2448621341701915; "P13" = 0.2448621341701916; "Q13" = 1071.900303233294; "R13" = 9647.102729099646; "S13" = 0.05391585235610947; "T13" = 0.05391585235610948
    "E14" = 3; "G14" = 7.798031999999999; "H14" = 23.394096; "I14" = 0.2874602628735535; "J14" = 0.2874602628735535; "K14" = 3; "M14" = 490.031855; "N14" = 1470.095565; "O14" = 0.6686419015677429; "P14" = 0.6686419015677431; "Q14" = 3821.28408630936; "R14" = 34391.55677678424; "S14" = 0.1922079767929361; "T14" = 0.1922079767929361
    "E15" = 3; "G15" = 7.798031999999999; "H15" = 23.394096; "I15" = 0.2874602628735535; "J15" = 0.2874602628735535; "K15" = 3; "M15" = 62.79306433333334; "N15" = 188.379193; "O15" = 0.0856802950924601; "P15" = 0.08568029509246011; "Q15" = 489.662325049392; "R15" = 4406.960925444528; "S15" = 0.02462968015036222; "T15" = 0.02462968015036222
    "E16" = 3; "G16" = 7.798031999999999; "H16" = 23.394096; "I16" = 0.2874602628735535; "J16" = 0.2874602628735535; "K16" = 3; "M16" = 0.5977846666666666; "N16" = 1.793354; "O16" = 0.0008156691696053909; "P16" = 0.000815669169605391; "Q16" = 4.661543959775999; "R16" = 41.95389563798399; "S16" = 0.0002344724739126188; "T16" = 0.0002344724739126188
    "E17" = 3; "G17" = 7.798031999999999; "H17" = 23.394096; "I17" = 0.2874602628735535; "J17" = 0.2874602628735535; "K17" = 3; "M17" = 179.453674; "N17" = 538.361022; "O17" = 0.2448621341701915; "P17" = 0.2448621341701916; "Q17" = 1399.385492369568; "R17" = 12594.46943132611; "S17" = 0.07038813345634258; "T17" = 0.07038813345634261
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
